# verigen : AddDefine 함수 추가
$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("summary")
$wsV = $wb.Worksheets.Item("_V")

# --- 1. Add new worksheet "AddDefine" after "_V" ---
$wsAdd = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsV)
$wsAdd.Name = "AddDefine"

# Column widths matching the "_V" sheet layout
$wsAdd.Columns.Item(1).ColumnWidth = 9.410714285714286
$wsAdd.Columns.Item(2).ColumnWidth = 55.160714285714285

# Page setup matching the other sheets
$wsAdd.PageSetup.PaperSize = 9
$wsAdd.PageSetup.Orientation = 1

# Row 1: function signature
$wsAdd.Range("A1").Value = "함수 원형"
$wsAdd.Range("A1").HorizontalAlignment = -4108
$wsAdd.Range("B1").Value = "function AddDefine(name, val)"

# Row 2: return value
$wsAdd.Range("A2").Value = "반환값"
$wsAdd.Range("A2").HorizontalAlignment = -4108
$wsAdd.Range("B2").Value = "-"

# Row 3: description (wrapped, taller row)
$wsAdd.Range("A3").Value = "설명"
$wsAdd.Range("A3").HorizontalAlignment = -4108
$wsAdd.Range("B3").Value = "코드내에서 적용될 전역 정의 변수를 지정합니다.`n아래의 함수 또는 지정된 code 내용에 적용됩니다.`nmodule:add_code, module.apply_code, module.code"
$wsAdd.Range("B3").WrapText = $true
$wsAdd.Rows.Item(3).RowHeight = 49.5

# Row 4: name parameter
$wsAdd.Range("A4").Value = "name"
$wsAdd.Range("A4").HorizontalAlignment = -4108
$wsAdd.Range("B4").Value = "정의 변수 이름 (재정의 가능)"

# Row 5: val parameter
$wsAdd.Range("A5").Value = "val"
$wsAdd.Range("A5").HorizontalAlignment = -4108
$wsAdd.Range("B5").Value = "정의 변수 값 (number 또는 string)"

$wsAdd.Range("B6").Select() | Out-Null

# --- 2. Update "summary" sheet with a new row referencing AddDefine ---
$wsSummary.Range("A3").Value = "AddDefine"
$wsSummary.Range("B3").Value = "function"
$wsSummary.Range("C3").Value = "전역 정의 선언"

$wsSummary.Range("C4").Select() | Out-Null

# --- 3. Update "_V" sheet selection ---
$wsV.Range("A1:B3").Select() | Out-Null

# --- 4. Leave "summary" as the active sheet/tab ---
$wsSummary.Activate()
$wsSummary.Range("C4").Select() | Out-Null
